# Update: Threat Alert Report - 2026-01-25 09:07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("THREAT_ALERT")

function Set-TextValue($range, $text) {
    # Assign as a formula that evaluates to the literal text, then freeze it
    # as a static value. This keeps the cell's existing style/number format
    # (e.g. "@"/General text) instead of letting Excel auto-detect a date
    # and silently reformat the cell with a new date number format/style.
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 2: 29-JAN-26 Nile Air NP-141 -> 05-FEB-26 Air Arabia Egypt E5-585 ---
Set-TextValue $ws.Range("A2") "05-FEB-26"
$ws.Range("C2").Value = "Air Arabia Egypt E5-585"
$ws.Range("D2").Value = 7184
$ws.Range("E2").Value = 7578
$ws.Range("F2").Value = -394

# --- Row 3: 31-JAN-26 Nile Air NP-141 -> 05-FEB-26 Air Arabia Egypt E5-585 ---
Set-TextValue $ws.Range("A3") "05-FEB-26"
$ws.Range("C3").Value = "Air Arabia Egypt E5-585"
$ws.Range("D3").Value = 7184
$ws.Range("E3").Value = 7531
$ws.Range("F3").Value = -347

# --- Row 4: 12-FEB-26 Air Arabia Egypt E5-585 -> 26-MAR-26 Nile Air NP-141 ---
Set-TextValue $ws.Range("A4") "26-MAR-26"
$ws.Range("C4").Value = "Nile Air NP-141"
$ws.Range("D4").Value = 12411
$ws.Range("E4").Value = 13512
$ws.Range("F4").Value = -1101

# --- Row 5: 26-FEB-26 -> 01-APR-26 ---
Set-TextValue $ws.Range("A5") "01-APR-26"
$ws.Range("D5").Value = 7725
$ws.Range("E5").Value = 8502
$ws.Range("F5").Value = -777

# --- Row 6: 15-APR-26 -> 13-MAY-26, threat upgraded to MEDIUM ---
Set-TextValue $ws.Range("A6") "13-MAY-26"
$ws.Range("D6").Value = 7725
$ws.Range("E6").Value = 10653
$ws.Range("F6").Value = -2928
$ws.Range("J10").Copy()
$ws.Range("J6").PasteSpecial(-4122)  # xlPasteFormats
Set-TextValue $ws.Range("J6") "MEDIUM THREAT - MONITOR"

# --- Row 7: 22-APR-26 -> 20-MAY-26, threat upgraded to MEDIUM ---
Set-TextValue $ws.Range("A7") "20-MAY-26"
$ws.Range("D7").Value = 7725
$ws.Range("E7").Value = 10653
$ws.Range("F7").Value = -2928
$ws.Range("J10").Copy()
$ws.Range("J7").PasteSpecial(-4122)  # xlPasteFormats
Set-TextValue $ws.Range("J7") "MEDIUM THREAT - MONITOR"

# --- Remove rows 8-12 (no longer part of the updated report) ---
$ws.Range("A8:K12").Delete()

$excel.CutCopyMode = $false
